$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 7773.143
$ws.Range("I28").Value = 8901.5
$ws.Range("J28").Value = 1003
$ws.Range("K28").Value = 8901.5
$ws.Range("L28").Value = 1003
$ws.Range("M28").Value = -8416.5
$ws.Range("N28").Value = -1973
$ws.Range("H53").Value = 58.42857
$ws.Range("I53").Value = 65.833336
$ws.Range("J53").Value = 14
$ws.Range("K53").Value = 65.833336
$ws.Range("L53").Value = 14
$ws.Range("M53").Value = 571.166664
$ws.Range("N53").Value = -1288
$ws.Range("H111").Value = 3465.25
$ws.Range("I111").Value = 3524.9285
$ws.Range("J111").Value = 3047.5
$ws.Range("K111").Value = 10574.7855
$ws.Range("L111").Value = 9142.5
$ws.Range("M111").Value = -7507.7855
$ws.Range("N111").Value = -15276.5
$ws.Range("H112").Value = 2234.087
$ws.Range("I112").Value = 1635.4286
$ws.Range("K112").Value = 4906.2858
$ws.Range("M112").Value = -3798.2858
$ws.Range("H115").Value = 1824.8334
$ws.Range("I115").Value = 1172.5454
$ws.Range("K115").Value = 3517.6362
$ws.Range("M115").Value = -1950.6362
$ws.Range("H116").Value = 2964.3333
$ws.Range("I116").Value = 2964.3333
$ws.Range("K116").Value = 2964.3333
$ws.Range("M116").Value = 477.6667000000002
$ws.Range("H132").Value = 1927.9354
$ws.Range("I132").Value = 1336.6666
$ws.Range("K132").Value = 4009.9998
$ws.Range("M132").Value = -1479.9998
$ws.Range("H141").Value = 4403
$ws.Range("I141").Value = 4214.778
$ws.Range("J141").Value = 5250
$ws.Range("K141").Value = 12644.334
$ws.Range("L141").Value = 15750
$ws.Range("M141").Value = -7464.334000000001
$ws.Range("N141").Value = -26110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5256.8076
$ws.Range("I32").Value = 5463.2
$ws.Range("J32").Value = 4568.8335
$ws.Range("K32").Value = 5463.2
$ws.Range("L32").Value = 4568.8335
$ws.Range("M32").Value = -5176.2
$ws.Range("N32").Value = -5142.8335
$ws.Range("H110").Value = 1846
$ws.Range("I110").Value = 710.5
$ws.Range("J110").Value = 2224.5
$ws.Range("K110").Value = 710.5
$ws.Range("L110").Value = 2224.5
$ws.Range("M110").Value = 1334.5
$ws.Range("N110").Value = -6314.5
$ws.Range("H113").Value = 60000
$ws.Range("J113").Value = 60000
$ws.Range("L113").Value = 60000
$ws.Range("N113").Value = -68678
$ws.Range("H124").Value = 52806.332
$ws.Range("J124").Value = 52806.332
$ws.Range("L124").Value = 52806.332
$ws.Range("N124").Value = -62626.332
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1268473.5
$ws.Range("I86").Value = 1925.75
$ws.Range("J86").Value = 2715956.8
$ws.Range("K86").Value = 1925.75
$ws.Range("L86").Value = 2715956.8
$ws.Range("M86").Value = -802.75
$ws.Range("N86").Value = -2718202.8
$ws.Range("H89").Value = 1268473.5
$ws.Range("I89").Value = 1925.75
$ws.Range("J89").Value = 2715956.8
$ws.Range("K89").Value = 9628.75
$ws.Range("L89").Value = 13579784
$ws.Range("M89").Value = -4012.75
$ws.Range("N89").Value = -13591016
$ws.Range("H105").Value = 1618.3572
$ws.Range("I105").Value = 1492.8462
$ws.Range("K105").Value = 1492.8462
$ws.Range("M105").Value = 254.1538
$ws.Range("H134").Value = 1693.463
$ws.Range("I134").Value = 1613.7826
$ws.Range("K134").Value = 4841.3478
$ws.Range("M134").Value = -2306.3478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2900
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 2900
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 2900
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -3180
$ws.Range("H31").Value = 2972.5
$ws.Range("I31").Value = 1312.1818
$ws.Range("J31").Value = 5255.4375
$ws.Range("K31").Value = 1312.1818
$ws.Range("L31").Value = 5255.4375
$ws.Range("M31").Value = -1017.1818
$ws.Range("N31").Value = -5845.4375
$ws.Range("H34").Value = 2972.5
$ws.Range("I34").Value = 1312.1818
$ws.Range("J34").Value = 5255.4375
$ws.Range("K34").Value = 1312.1818
$ws.Range("L34").Value = 5255.4375
$ws.Range("M34").Value = -1110.1818
$ws.Range("N34").Value = -5659.4375
$ws.Range("H99").Value = 3392.7144
$ws.Range("J99").Value = 3766.6667
$ws.Range("L99").Value = 3766.6667
$ws.Range("N99").Value = -6762.6667
$ws.Range("H126").Value = 3392.7144
$ws.Range("J126").Value = 3766.6667
$ws.Range("L126").Value = 11300.0001
$ws.Range("N126").Value = -16240.0001
$ws.Range("H132").Value = 1553.3572
$ws.Range("I132").Value = 1462.3636
$ws.Range("J132").Value = 1887
$ws.Range("K132").Value = 4387.0908
$ws.Range("L132").Value = 5661
$ws.Range("M132").Value = -1857.0908
$ws.Range("N132").Value = -10721
$ws.Range("H134").Value = 2079.5625
$ws.Range("I134").Value = 2051.6
$ws.Range("K134").Value = 6154.799999999999
$ws.Range("M134").Value = -3619.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7470451.5
$ws.Range("I4").Value = 1777600.6
$ws.Range("J4").Value = 50166830
$ws.Range("K4").Value = 5332801.800000001
$ws.Range("L4").Value = 150500490
$ws.Range("M4").Value = -5332689.800000001
$ws.Range("N4").Value = -150500714
$ws.Range("H80").Value = 1798.5
$ws.Range("I80").Value = 1798.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 5395.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -4459.5
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 1798.5
$ws.Range("I83").Value = 1798.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 16186.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -11506.5
$ws.Range("N83").Value = ""
$ws.Range("H131").Value = 1687.7778
$ws.Range("I131").Value = 1057.091
$ws.Range("J131").Value = 1965.28
$ws.Range("K131").Value = 3171.273
$ws.Range("L131").Value = 5895.84
$ws.Range("M131").Value = 1868.727
$ws.Range("N131").Value = -15975.84
$ws.Range("H132").Value = 1488.4348
$ws.Range("I132").Value = 802.4286
$ws.Range("K132").Value = 7221.8574
$ws.Range("M132").Value = -4691.8574
$ws.Range("H138").Value = 4900
$ws.Range("I138").Value = 3250
$ws.Range("J138").Value = 6000
$ws.Range("K138").Value = 9750
$ws.Range("L138").Value = 18000
$ws.Range("M138").Value = -4610
$ws.Range("N138").Value = -28280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23771.4
$ws.Range("J15").Value = 23771.4
$ws.Range("L15").Value = 23771.4
$ws.Range("N15").Value = -24347.4
$ws.Range("H81").Value = 23771.4
$ws.Range("J81").Value = 23771.4
$ws.Range("L81").Value = 23771.4
$ws.Range("N81").Value = -25767.4
$ws.Range("H84").Value = 23771.4
$ws.Range("J84").Value = 23771.4
$ws.Range("L84").Value = 71314.20000000001
$ws.Range("N84").Value = -81298.20000000001
$ws.Range("H107").Value = 1987.8667
$ws.Range("I107").Value = 227.375
$ws.Range("K107").Value = 227.375
$ws.Range("M107").Value = 1692.625
$ws.Range("H122").Value = 2978.1904
$ws.Range("I122").Value = 2292.8333
$ws.Range("J122").Value = 3892
$ws.Range("K122").Value = 6878.499899999999
$ws.Range("L122").Value = 11676
$ws.Range("M122").Value = -4428.499899999999
$ws.Range("N122").Value = -16576
$ws.Range("H123").Value = 23998.8
$ws.Range("J123").Value = 23998.8
$ws.Range("L123").Value = 23998.8
$ws.Range("N123").Value = -28898.8
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 3500
$ws.Range("K126").Value = 10500
$ws.Range("M126").Value = -8030
$ws.Range("H132").Value = 2122.2727
$ws.Range("I132").Value = 1893.7368
$ws.Range("K132").Value = 5681.2104
$ws.Range("M132").Value = -3151.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2846.6287
$ws.Range("I22").Value = 2322.111
$ws.Range("K22").Value = 2322.111
$ws.Range("M22").Value = -2027.111
$ws.Range("H27").Value = 2846.6287
$ws.Range("I27").Value = 2322.111
$ws.Range("K27").Value = 2322.111
$ws.Range("M27").Value = -2215.111
$ws.Range("H61").Value = 3552.1667
$ws.Range("I61").Value = 2528.55
$ws.Range("K61").Value = 2528.55
$ws.Range("M61").Value = -2326.55
$ws.Range("H100").Value = 8000.1
$ws.Range("I100").Value = 5500.25
$ws.Range("K100").Value = 5500.25
$ws.Range("M100").Value = -4959.25
$ws.Range("H113").Value = 3552.1667
$ws.Range("I113").Value = 2528.55
$ws.Range("K113").Value = 2528.55
$ws.Range("M113").Value = -358.5500000000002
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
